# Case with 380 kV done: update the slack/external-grid setpoint (B column, 1.05 -> 1.02 pu)
# and refresh the resulting power-flow bus voltage magnitudes (res_bus/vm_pu.xlsx) for every
# bus row (rows 2-25). Column G (=1) and the empty H column are not affected by the rerun.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.026250579526633
$ws.Range("D2").Value = 1.053931840132151
$ws.Range("E2").Value = 1.038910667876532
$ws.Range("F2").Value = 1.055748805978192
$ws.Range("I2").Value = 1.04347076485328
$ws.Range("J2").Value = 1.031415066543644
$ws.Range("K2").Value = 1.056676446262077
$ws.Range("L2").Value = 1.041697355734919
$ws.Range("M2").Value = 1.058488409858822
$ws.Range("N2").Value = 1.014500921118896
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.027158924581206
$ws.Range("D3").Value = 1.054581360211444
$ws.Range("E3").Value = 1.039691773165726
$ws.Range("F3").Value = 1.056553533799975
$ws.Range("I3").Value = 1.043629259391321
$ws.Range("J3").Value = 1.031963303574436
$ws.Range("K3").Value = 1.057139527880136
$ws.Range("L3").Value = 1.042288647969613
$ws.Range("M3").Value = 1.059106662827597
$ws.Range("N3").Value = 1.014683003432843
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.027747129478519
$ws.Range("D4").Value = 1.055000095123257
$ws.Range("E4").Value = 1.040197480597366
$ws.Range("F4").Value = 1.057073463683359
$ws.Range("I4").Value = 1.043729419100503
$ws.Range("J4").Value = 1.032317879723355
$ws.Range("K4").Value = 1.057436941424941
$ws.Range("L4").Value = 1.04267087643754
$ws.Range("M4").Value = 1.059505278445231
$ws.Range("N4").Value = 1.014800734056865
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.027994516195931
$ws.Range("D5").Value = 1.055175758548427
$ws.Range("E5").Value = 1.040410145443765
$ws.Range("F5").Value = 1.057291852834006
$ws.Range("I5").Value = 1.043770951292097
$ws.Range("J5").Value = 1.032466901831127
$ws.Range("K5").Value = 1.057561438073082
$ws.Range("L5").Value = 1.042831473686673
$ws.Range("M5").Value = 1.059672511043899
$ws.Range("N5").Value = 1.014850206375642
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.028036059660763
$ws.Range("D6").Value = 1.055205231290668
$ws.Range("E6").Value = 1.040445856597764
$ws.Range("F6").Value = 1.057328510170977
$ws.Range("I6").Value = 1.043777890981019
$ws.Range("J6").Value = 1.032491920822803
$ws.Range("K6").Value = 1.057582310122257
$ws.Range("L6").Value = 1.042858433261764
$ws.Range("M6").Value = 1.059700569826308
$ws.Range("N6").Value = 1.014858511715992
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.027750434659311
$ws.Range("D7").Value = 1.055002443813877
$ws.Range("E7").Value = 1.040200321980361
$ws.Range("F7").Value = 1.057076382556506
$ws.Range("I7").Value = 1.043729976316313
$ws.Range("J7").Value = 1.032319871129247
$ws.Range("K7").Value = 1.057438607063772
$ws.Range("L7").Value = 1.042673022707438
$ws.Range("M7").Value = 1.059507514374895
$ws.Range("N7").Value = 1.014801395193899
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.026557466189795
$ws.Range("D8").Value = 1.05415166773525
$ws.Range("E8").Value = 1.039174587110455
$ws.Range("F8").Value = 1.0560209284306
$ws.Range("I8").Value = 1.043524824430262
$ws.Range("J8").Value = 1.031600380365927
$ws.Range("K8").Value = 1.056833407751942
$ws.Range("L8").Value = 1.041897263021386
$ws.Range("M8").Value = 1.058697647494623
$ws.Range("N8").Value = 1.01456247485339
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.024458760416448
$ws.Range("D9").Value = 1.052640728521675
$ws.Range("E9").Value = 1.037369326685384
$ws.Range("F9").Value = 1.054155173306565
$ws.Range("I9").Value = 1.04314501290252
$ws.Range("J9").Value = 1.030331290105887
$ws.Range("K9").Value = 1.055749960845857
$ws.Range("L9").Value = 1.040527444004273
$ws.Range("M9").Value = 1.05725964217173
$ws.Range("N9").Value = 1.014140802916194
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.023062012806506
$ws.Range("D10").Value = 1.051625649744886
$ws.Range("E10").Value = 1.036167402775326
$ws.Range("F10").Value = 1.052907475788311
$ws.Range("I10").Value = 1.042879561204314
$ws.Range("J10").Value = 1.029484445757916
$ws.Range("K10").Value = 1.055016333205616
$ws.Range("L10").Value = 1.039612403345939
$ws.Range("M10").Value = 1.056293733996571
$ws.Range("N10").Value = 1.013859263747251
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.022457785669508
$ws.Range("D11").Value = 1.051184287912845
$ws.Range("E11").Value = 1.035647351216967
$ws.Range("F11").Value = 1.052366316369525
$ws.Range("I11").Value = 1.042761726828559
$ws.Range("J11").Value = 1.029117579109925
$ws.Range("K11").Value = 1.054695998708712
$ws.Range("L11").Value = 1.039215762592557
$ws.Range("M11").Value = 1.055873789111227
$ws.Range("N11").Value = 1.013737257953926
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.022233436022771
$ws.Range("D12").Value = 1.051020074096409
$ws.Range("E12").Value = 1.035454240936622
$ws.Range("F12").Value = 1.052165172519576
$ws.Range("I12").Value = 1.042717524403581
$ws.Range("J12").Value = 1.028981282793759
$ws.Range("K12").Value = 1.054576612906694
$ws.Range("L12").Value = 1.039068370345669
$ws.Range("M12").Value = 1.055717549036313
$ws.Range("N12").Value = 1.013691925273256
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.022281555820741
$ws.Range("D13").Value = 1.051055310812084
$ws.Range("E13").Value = 1.035495661006859
$ws.Range("F13").Value = 1.052208324554817
$ws.Range("I13").Value = 1.042727025574215
$ws.Range("J13").Value = 1.029010519964385
$ws.Range("K13").Value = 1.054602239609532
$ws.Range("L13").Value = 1.03909998927887
$ws.Range("M13").Value = 1.055751074524368
$ws.Range("N13").Value = 1.013701649928542
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.022439239072738
$ws.Range("D14").Value = 1.051170719492799
$ws.Range("E14").Value = 1.035631387427761
$ws.Range("F14").Value = 1.052349692462012
$ws.Range("I14").Value = 1.042758081879659
$ws.Range("J14").Value = 1.029106313332874
$ws.Range("K14").Value = 1.05468613837729
$ws.Range("L14").Value = 1.039203580365592
$ws.Range("M14").Value = 1.05586087942711
$ws.Range("N14").Value = 1.013733511030046
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.022536404515468
$ws.Range("D15").Value = 1.051241790555038
$ws.Range("E15").Value = 1.035715020954098
$ws.Range("F15").Value = 1.052436776305348
$ws.Range("I15").Value = 1.042777159282408
$ws.Range("J15").Value = 1.02916533141054
$ws.Range("K15").Value = 1.054737778314374
$ws.Range("L15").Value = 1.039267398045436
$ws.Range("M15").Value = 1.055928500262081
$ws.Range("N15").Value = 1.01375313982214
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.023102125547027
$ws.Range("D16").Value = 1.051654903171713
$ws.Range("E16").Value = 1.036201925219029
$ws.Range("F16").Value = 1.052943371986913
$ws.Range("I16").Value = 1.042887320637973
$ws.Range("J16").Value = 1.029508789831755
$ws.Range("K16").Value = 1.055037536663339
$ws.Range("L16").Value = 1.039638718296323
$ws.Range("M16").Value = 1.056321568669866
$ws.Range("N16").Value = 1.013867358845992
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.023457141738811
$ws.Range("D17").Value = 1.051913550315238
$ws.Range("E17").Value = 1.036507452735864
$ws.Range("F17").Value = 1.053260906882918
$ws.Range("I17").Value = 1.042955648379697
$ws.Range("J17").Value = 1.02972418538987
$ws.Range("K17").Value = 1.055224853801516
$ws.Range("L17").Value = 1.039871525524273
$ws.Range("M17").Value = 1.056567676022936
$ws.Range("N17").Value = 1.01393897957148
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.023664271972576
$ws.Range("D18").Value = 1.052064238503474
$ws.Range("E18").Value = 1.036685699240278
$ws.Range("F18").Value = 1.053446032837172
$ws.Range("I18").Value = 1.042995223689159
$ws.Range("J18").Value = 1.02984980478922
$ws.Range("K18").Value = 1.055333855100001
$ws.Range("L18").Value = 1.040007277233791
$ws.Range("M18").Value = 1.056711062120527
$ws.Range("N18").Value = 1.01398074530079
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.023734907377885
$ws.Range("D19").Value = 1.052115589315512
$ws.Range("E19").Value = 1.036746483001711
$ws.Range("F19").Value = 1.053509141229546
$ws.Range("I19").Value = 1.043008670474888
$ws.Range("J19").Value = 1.029892634814056
$ws.Range("K19").Value = 1.0553709779783
$ws.Range("L19").Value = 1.040053558100013
$ws.Range("M19").Value = 1.056759925181392
$ws.Range("N19").Value = 1.013994984725366
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.023419046134444
$ws.Range("D20").Value = 1.051885818147319
$ws.Range("E20").Value = 1.036474668637761
$ws.Range("F20").Value = 1.053226847356523
$ws.Range("I20").Value = 1.042948346321586
$ws.Range("J20").Value = 1.029701077255609
$ws.Range("K20").Value = 1.055204783079948
$ws.Range("L20").Value = 1.039846551724132
$ws.Range("M20").Value = 1.056541288001008
$ws.Range("N20").Value = 1.013931296318114
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.022392802869527
$ws.Range("D21").Value = 1.051136742016247
$ws.Range("E21").Value = 1.035591417710628
$ws.Range("F21").Value = 1.052308066798065
$ws.Range("I21").Value = 1.042748948530255
$ws.Range("J21").Value = 1.029078105275503
$ws.Range("K21").Value = 1.054661443290903
$ws.Range("L21").Value = 1.039173077080178
$ws.Range("M21").Value = 1.055828551614728
$ws.Range("N21").Value = 1.013724129121649
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.021748068094211
$ws.Range("D22").Value = 1.050664193163342
$ws.Range("E22").Value = 1.035036430817327
$ws.Range("F22").Value = 1.05172962352018
$ws.Range("I22").Value = 1.042621071663779
$ws.Range("J22").Value = 1.028686269791628
$ws.Range("K22").Value = 1.054317514538488
$ws.Range("L22").Value = 1.038749277498777
$ws.Range("M22").Value = 1.055378958059678
$ws.Range("N22").Value = 1.013593792336807
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.022089805984159
$ws.Range("D23").Value = 1.050914848892045
$ws.Range("E23").Value = 1.035330606378598
$ws.Range("F23").Value = 1.052036339591095
$ws.Range("I23").Value = 1.042689098966789
$ws.Range("J23").Value = 1.028894002933928
$ws.Range("K23").Value = 1.05450005603913
$ws.Range("L23").Value = 1.038973975337096
$ws.Range("M23").Value = 1.055617434721163
$ws.Range("N23").Value = 1.013662894041406
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.023436259727446
$ws.Range("D24").Value = 1.051898349665207
$ws.Range("E24").Value = 1.036489482244861
$ws.Range("F24").Value = 1.053242237658676
$ws.Range("I24").Value = 1.04295164666978
$ws.Range("J24").Value = 1.029711518879969
$ws.Range("K24").Value = 1.055213852971193
$ws.Range("L24").Value = 1.039857836435606
$ws.Range("M24").Value = 1.056553212119576
$ws.Range("N24").Value = 1.013934768078746
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.025000909743112
$ws.Range("D25").Value = 1.053032723742363
$ws.Range("E25").Value = 1.037835757550234
$ws.Range("F25").Value = 1.05463820438907
$ws.Range("I25").Value = 1.043245366073698
$ws.Range("J25").Value = 1.03065952293909
$ws.Range("K25").Value = 1.056032062421605
$ws.Range("L25").Value = 1.040881902740547
$ws.Range("M25").Value = 1.057632684109652
$ws.Range("N25").Value = 1.014249891629766
